$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").ClearContents()
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.186.56"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").ClearContents()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.86"
$ws.Range("E3").Value = "  +3.11%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").ClearContents()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.83"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").ClearContents()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5037"
$ws.Range("E7").Value = "  -1.29%  "

$ws.Range("D8").ClearContents()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3924"
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").ClearContents()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09587"
$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").ClearContents()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  +3.29%  "

$ws.Range("D11").ClearContents()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.91"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").ClearContents()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.497"
$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("D13").ClearContents()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +2.98%  "

$ws.Range("D14").ClearContents()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.870.60"
$ws.Range("E14").Value = "  +3.36%  "

$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").ClearContents()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.419"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("D17").ClearContents()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001131"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("D18").ClearContents()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.05"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").ClearContents()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06627"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").ClearContents()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.53"
$ws.Range("E20").Value = "  +1.68%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").ClearContents()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.159"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("D23").ClearContents()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.258.64"
$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("D24").ClearContents()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").ClearContents()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  +2.81%  "

$ws.Range("D26").ClearContents()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.529"
$ws.Range("E26").Value = "  +4.37%  "

$ws.Range("D27").ClearContents()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.085.41"
$ws.Range("E27").Value = "  +3.77%  "

$ws.Range("D28").ClearContents()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.19"
$ws.Range("E28").Value = "  +3.52%  "

$ws.Range("D29").ClearContents()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.70"
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("D30").ClearContents()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.56"
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("D31").ClearContents()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.068"
$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").ClearContents()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1056"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").ClearContents()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.624"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").ClearContents()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.627"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").ClearContents()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.607"
$ws.Range("E35").Value = "  +6.05%  "

$ws.Range("D36").ClearContents()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06752"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").ClearContents()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02387"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").ClearContents()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2186"
$ws.Range("E38").Value = "  +0.93%  "

$ws.Range("D39").ClearContents()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "11.46"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").ClearContents()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6342"
$ws.Range("E40").Value = "  +2.84%  "

$ws.Range("D41").ClearContents()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.978"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("E42").Value = "  +2.31%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").ClearContents()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  +1.95%  "

$ws.Range("D45").ClearContents()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6024"
$ws.Range("E45").Value = "  +1.44%  "

$ws.Range("D46").ClearContents()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.667"
$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").ClearContents()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.267"
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").ClearContents()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.04"
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("D49").ClearContents()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").ClearContents()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").ClearContents()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06840"
$ws.Range("E51").Value = "  +1.12%  "
